$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C holds a "Förändrad" (changed) date as a serial number.
# Every row from 2 to 385 had its value bumped from 45179 to 45180
# (one day later).
$ws.Range("C2:C385").Value = 45180
